# Release refresh: bump the embedded build timestamp in the version string
# from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet (A2, A6) and every row of the "S" column on the
# "Boundaries and methane sources" sheet (S2:S12).

$wb = $excel.ActiveWorkbook

$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

$aboutWs.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Zhaozhuang Coal Mine, China, M0430, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

for ($row = 2; $row -le 12; $row++) {
    $dataWs.Range("S" + $row).Value = $newVersion
}
